$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New wallet-fund rows (S.No 59-63), landing on sheet rows 62-66.
# Columns: A=S.No B=Order ID C=Topup Amount D=Order Amount E=Creation Date
$data = @(
    @(59, 26612971, 47515, 45000,    44286),
    @(60, 26645301, 68583, 64999.6,  44287),
    @(61, 26657499, 63308, 60000.1,  44288),
    @(62, 26378669, 52757, 50000.08, 44289),
    @(63, 26666554, 79135, 75000.62, 44290)
)

$row = 62
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value2 = $entry[0]
    $ws.Cells.Item($row, 2).Value2 = $entry[1]
    $ws.Cells.Item($row, 3).Value2 = $entry[2]
    $ws.Cells.Item($row, 4).Value2 = $entry[3]
    $ws.Cells.Item($row, 5).Value2 = $entry[4]

    # Match the formatting of column A used by the rows above (col B keeps
    # the plain default format it already received).
    $ws.Range("A61").Copy() | Out-Null
    $ws.Range("A" + $row).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    # Re-assert the existing shared formulas so the cached values recalc.
    $ws.Cells.Item($row, 6).Formula = '=IF(B' + $row + '="","",C' + $row + '-D' + $row + ')'
    $ws.Cells.Item($row, 7).Formula = '=IF(B' + $row + '="","",F' + $row + '/D' + $row + '*100)'
    $ws.Cells.Item($row, 8).Formula = '=IF(B' + $row + '="","",D' + $row + '*1.04)'
    $ws.Cells.Item($row, 9).Formula = '=IF(B' + $row + '="","",C' + $row + '-H' + $row + ')'

    $row = $row + 1
}

$excel.Calculate()

$ws.Range("C64:D64").Select() | Out-Null

$wb.Save()
